# Auto-generated edit script applying the diff to Cerberus_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2618.6
$ws.Cells.Item(40, 9).Value = 2198
$ws.Cells.Item(40, 10).Value = 3249.5
$ws.Cells.Item(40, 11).Value = 2198
$ws.Cells.Item(40, 12).Value = 3249.5
$ws.Cells.Item(40, 13).Value = -2023
$ws.Cells.Item(40, 14).Value = -3599.5

$ws.Cells.Item(62, 8).Value = 25002566
$ws.Cells.Item(62, 9).Value = 27780218
$ws.Cells.Item(62, 11).Value = 27780218
$ws.Cells.Item(62, 13).Value = -27779594

$ws.Cells.Item(65, 8).Value = 25002566
$ws.Cells.Item(65, 9).Value = 27780218
$ws.Cells.Item(65, 11).Value = 138901090
$ws.Cells.Item(65, 13).Value = -138897970

$ws.Cells.Item(70, 8).Value = 14974
$ws.Cells.Item(70, 9).Value = 4865
$ws.Cells.Item(70, 11).Value = 14595
$ws.Cells.Item(70, 13).Value = -14325

$ws.Cells.Item(73, 8).Value = 14974
$ws.Cells.Item(73, 9).Value = 4865
$ws.Cells.Item(73, 11).Value = 14595
$ws.Cells.Item(73, 13).Value = -13659

$ws.Cells.Item(88, 8).Value = 1871.75
$ws.Cells.Item(88, 9).Value = 1857.6
$ws.Cells.Item(88, 11).Value = 1857.6
$ws.Cells.Item(88, 13).Value = -1451.6

$ws.Cells.Item(91, 8).Value = 1871.75
$ws.Cells.Item(91, 9).Value = 1857.6
$ws.Cells.Item(91, 11).Value = 1857.6
$ws.Cells.Item(91, 13).Value = -453.5999999999999

$ws.Cells.Item(96, 8).Value = 690
$ws.Cells.Item(96, 9).Value = 649.6667
$ws.Cells.Item(96, 10).Value = 738.4
$ws.Cells.Item(96, 11).Value = 1949.0001
$ws.Cells.Item(96, 12).Value = 2215.2
$ws.Cells.Item(96, 13).Value = -576.0001
$ws.Cells.Item(96, 14).Value = -4961.2

$ws.Cells.Item(106, 8).Value = 6514712.5
$ws.Cells.Item(106, 9).Value = 8286016.5
$ws.Cells.Item(106, 10).Value = 19932.666
$ws.Cells.Item(106, 11).Value = 8286016.5
$ws.Cells.Item(106, 12).Value = 19932.666
$ws.Cells.Item(106, 13).Value = -8285385.5
$ws.Cells.Item(106, 14).Value = -21194.666

$ws.Cells.Item(116, 8).Value = 10427.789
$ws.Cells.Item(116, 9).Value = 12628.143
$ws.Cells.Item(116, 10).Value = 9144.25
$ws.Cells.Item(116, 11).Value = 12628.143
$ws.Cells.Item(116, 12).Value = 9144.25
$ws.Cells.Item(116, 13).Value = -9186.143
$ws.Cells.Item(116, 14).Value = -16028.25

$ws.Cells.Item(137, 8).Value = 835502
$ws.Cells.Item(137, 9).Value = 1667915.1
$ws.Cells.Item(137, 10).Value = 3088.8333
$ws.Cells.Item(137, 11).Value = 5003745.300000001
$ws.Cells.Item(137, 12).Value = 9266.499899999999
$ws.Cells.Item(137, 13).Value = -5001195.300000001
$ws.Cells.Item(137, 14).Value = -14366.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1973.25
$ws.Cells.Item(63, 9).Value = 1973.25
$ws.Cells.Item(63, 11).Value = 1973.25
$ws.Cells.Item(63, 13).Value = -1287.25

$ws.Cells.Item(66, 8).Value = 1973.25
$ws.Cells.Item(66, 9).Value = 1973.25
$ws.Cells.Item(66, 11).Value = 9866.25
$ws.Cells.Item(66, 13).Value = -6434.25

$ws.Cells.Item(74, 8).Value = 2166.7878
$ws.Cells.Item(74, 9).Value = 1200.6471
$ws.Cells.Item(74, 10).Value = 3193.3125
$ws.Cells.Item(74, 11).Value = 1200.6471
$ws.Cells.Item(74, 12).Value = 3193.3125
$ws.Cells.Item(74, 13).Value = -326.6470999999999
$ws.Cells.Item(74, 14).Value = -4941.3125

$ws.Cells.Item(77, 8).Value = 2166.7878
$ws.Cells.Item(77, 9).Value = 1200.6471
$ws.Cells.Item(77, 10).Value = 3193.3125
$ws.Cells.Item(77, 11).Value = 6003.2355
$ws.Cells.Item(77, 12).Value = 15966.5625
$ws.Cells.Item(77, 13).Value = -1635.2355
$ws.Cells.Item(77, 14).Value = -24702.5625

$ws.Cells.Item(122, 8).Value = 1683.2941
$ws.Cells.Item(122, 9).Value = 1724.32
$ws.Cells.Item(122, 10).Value = 1569.3334
$ws.Cells.Item(122, 11).Value = 5172.96
$ws.Cells.Item(122, 12).Value = 4708.0002
$ws.Cells.Item(122, 13).Value = -2722.96
$ws.Cells.Item(122, 14).Value = -9608.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3096.8696
$ws.Cells.Item(31, 9).Value = 2045.6875
$ws.Cells.Item(31, 11).Value = 2045.6875
$ws.Cells.Item(31, 13).Value = -1750.6875

$ws.Cells.Item(34, 8).Value = 3096.8696
$ws.Cells.Item(34, 9).Value = 2045.6875
$ws.Cells.Item(34, 11).Value = 2045.6875
$ws.Cells.Item(34, 13).Value = -1843.6875

$ws.Cells.Item(105, 8).Value = 1307.1428
$ws.Cells.Item(105, 9).Value = 1200.4117
$ws.Cells.Item(105, 11).Value = 1200.4117
$ws.Cells.Item(105, 13).Value = 546.5882999999999

$ws.Cells.Item(134, 8).Value = 10210.3
$ws.Cells.Item(134, 9).Value = 9099
$ws.Cells.Item(134, 11).Value = 27297
$ws.Cells.Item(134, 13).Value = -24762

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 30665322
$ws.Cells.Item(4, 9).Value = 45950644
$ws.Cells.Item(4, 10).Value = 94681.55
$ws.Cells.Item(4, 11).Value = 137851932
$ws.Cells.Item(4, 12).Value = 284044.65
$ws.Cells.Item(4, 13).Value = -137851820
$ws.Cells.Item(4, 14).Value = -284268.65

$ws.Cells.Item(41, 8).Value = 9917.615
$ws.Cells.Item(41, 10).Value = 14082.4
$ws.Cells.Item(41, 12).Value = 42247.2
$ws.Cells.Item(41, 14).Value = -42923.2

$ws.Cells.Item(80, 8).Value = 5666
$ws.Cells.Item(80, 10).Value = 5666
$ws.Cells.Item(80, 12).Value = 16998
$ws.Cells.Item(80, 14).Value = -18870

$ws.Cells.Item(83, 8).Value = 5666
$ws.Cells.Item(83, 10).Value = 5666
$ws.Cells.Item(83, 12).Value = 50994
$ws.Cells.Item(83, 14).Value = -60354

$ws.Cells.Item(98, 8).Value = 41671140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 6750.25
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 2622.0557
$ws.Cells.Item(80, 10).Value = 2380.625
$ws.Cells.Item(80, 12).Value = 2380.625
$ws.Cells.Item(80, 14).Value = -4376.625

$ws.Cells.Item(83, 8).Value = 2622.0557
$ws.Cells.Item(83, 10).Value = 2380.625
$ws.Cells.Item(83, 12).Value = 11903.125
$ws.Cells.Item(83, 14).Value = -21887.125

$ws.Cells.Item(97, 8).Value = 871.55554
$ws.Cells.Item(97, 9).Value = 372.6
$ws.Cells.Item(97, 11).Value = 372.6
$ws.Cells.Item(97, 13).Value = 123.4

$ws.Cells.Item(102, 8).Value = 4151.7666
$ws.Cells.Item(102, 9).Value = 4205.5
$ws.Cells.Item(102, 10).Value = 3399.5
$ws.Cells.Item(102, 11).Value = 4205.5
$ws.Cells.Item(102, 12).Value = 3399.5
$ws.Cells.Item(102, 13).Value = -2583.5
$ws.Cells.Item(102, 14).Value = -6643.5

$ws.Cells.Item(132, 8).Value = 4664.6665
$ws.Cells.Item(132, 9).Value = 4711.375
$ws.Cells.Item(132, 11).Value = 14134.125
$ws.Cells.Item(132, 13).Value = -11604.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 3999
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(22, 8).Value = 1362.3334
$ws.Cells.Item(22, 9).Value = 837
$ws.Cells.Item(22, 11).Value = 837
$ws.Cells.Item(22, 13).Value = -542

$ws.Cells.Item(27, 8).Value = 1362.3334
$ws.Cells.Item(27, 9).Value = 837
$ws.Cells.Item(27, 11).Value = 837
$ws.Cells.Item(27, 13).Value = -730

$ws.Cells.Item(82, 8).Value = 3626.7778
$ws.Cells.Item(82, 10).Value = 7839
$ws.Cells.Item(82, 12).Value = 7839
$ws.Cells.Item(82, 14).Value = -8561

$ws.Cells.Item(85, 8).Value = 3626.7778
$ws.Cells.Item(85, 10).Value = 7839
$ws.Cells.Item(85, 12).Value = 7839
$ws.Cells.Item(85, 14).Value = -10335

$ws.Cells.Item(132, 8).Value = 3970
$ws.Cells.Item(132, 9).Value = 3762.5
$ws.Cells.Item(132, 11).Value = 11287.5
$ws.Cells.Item(132, 13).Value = -8757.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8499
$ws.Cells.Item(62, 9).Value = 6999.5
$ws.Cells.Item(62, 11).Value = 6999.5
$ws.Cells.Item(62, 13).Value = -6375.5

$ws.Cells.Item(65, 8).Value = 8499
$ws.Cells.Item(65, 9).Value = 6999.5
$ws.Cells.Item(65, 11).Value = 34997.5
$ws.Cells.Item(65, 13).Value = -31877.5

$ws.Cells.Item(132, 8).Value = 2500
$ws.Cells.Item(132, 9).Value = 2500
$ws.Cells.Item(132, 11).Value = 7500
$ws.Cells.Item(132, 13).Value = -4970
